$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.612.70'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '2.339.86'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.22'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.73%  '
$ws.Range('E7').Value = '  -2.12%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.517'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.34'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.67%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '51.82'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.43%  '
$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0798'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.94%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.113'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.83'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.710.12'
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.56'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.348.14'
$ws.Range('E17').Value = '  -0.97%  '
$ws.Range('B18').Value = 'Polygon'
$ws.Range('C18').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.808'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.98%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '43.531.42'
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.56%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0908'
$ws.Range('E21').Value = '  -1.97%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.13'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.33%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.19'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.46%  '
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '238.28'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.52%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.99'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.32%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.55'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.12%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.03'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.18'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.85%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.67'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '166.17'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.61%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.28'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.58%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.07'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.28%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.42'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.82%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.50'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.36%  '
$ws.Range('B37').Value = 'Celestia'
$ws.Range('C37').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '16.95'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.56%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0708'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.65%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.91'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.24%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.83'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.10%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.103'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.73%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.111'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.75%  '
$ws.Range('D43').Value = '1.992.13'
$ws.Range('E43').Value = '  -0.49%  '
$ws.Range('E44').Value = '  -1.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.56'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.63%  '
$ws.Range('E46').Value = '  -7.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '56.71'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.71%  '
$ws.Range('E49').Value = '  +4.06%  '
$ws.Range('D50').Value = '2.567.05'
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('E51').Value = '  -1.09%  '
